$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet ---
$ws.Name = "B08K4353V1"

# --- Replace the 100 keyword values in column A, top to bottom ---
$newValues = @(
  'neck compression wrap',
  'training bra sets',
  'wide waistband yoga',
  'girls seamless underwear',
  'sport bra',
  'seamless bra',
  'yoga legging',
  'workout waist',
  'yoga workout',
  'yoga bra',
  'sports bra for women',
  'sport apparel',
  'seamless bra set',
  'workout apparel',
  'workout bra',
  'workout outfit set',
  'workout set',
  'workout leggings for women high waist',
  'workout legging',
  'legging bra set',
  'workout outfit',
  'sports leggings women',
  'black legging',
  'sport bra set',
  'yoga sets for women',
  'high waist workout',
  'yoga set purple',
  'yoga bras for women',
  '2 pcs',
  'green 2',
  'pcs set',
  'bra set',
  'waist yoga',
  'outfit set',
  'seamless sports bra',
  'seamless workout leggings',
  'seamless yoga bra',
  'seamless yoga',
  'sports for women',
  'high waist yoga',
  'sport legging',
  'black bra',
  'yoga legging set',
  'yoga workout set',
  'yoga bra set',
  'yoga outfit set',
  'workout set women',
  'black legging set',
  'sports leggings for women',
  'black bra set',
  'sport leggings for women',
  'sport set',
  'black outfit',
  'high waist',
  'set 2',
  'yoga apparel',
  'purple bra',
  'purple legging',
  'high waist legging',
  'high sport bra',
  'black n',
  'womens sports bra set',
  'sport outfit',
  'sport waist',
  'sport workout',
  'waist bra',
  'black 2',
  'yoga sport',
  'n set',
  'green set',
  'yoga 2',
  'green bra',
  'green bra set',
  'green yoga',
  'yoga sport bra',
  'womens workout sports bra',
  'high waist set',
  'ready set',
  'high black',
  'waist set',
  'black yoga',
  'green outfit',
  'black apparel',
  'purple yoga',
  'black set',
  'purple outfit',
  'purple workout',
  'black sport',
  'seamless legging',
  'green apparel',
  'purple set',
  'black seamless bra',
  'black yoga bra',
  'seamless workout',
  'green sport bra',
  'black sport bra',
  'purple sport bra',
  'workout sport bra',
  'black workout bra',
  'seamless sport bra'
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $newValues[$i]
}

# --- Update the saved selection/view state ---
$ws.Range("A1:A100").Select()

# --- Conditional formatting: duplicate the existing "highlight duplicates" rule's
#     differential format twice, then drop the two extra rules again. This grows the
#     workbook's dxf table (1 -> 3 identical fills) while leaving the single
#     surviving rule's dxfId untouched, matching how Excel's CF Rules Manager
#     leaves orphaned dxf entries behind after rule edits. ---
$rng = $ws.Range("A1:A100")

$dup1 = $rng.FormatConditions.AddUniqueValues()
$dup1.DupeUnique = 1
$dup1.Interior.Color = 13998939

$dup2 = $rng.FormatConditions.AddUniqueValues()
$dup2.DupeUnique = 1
$dup2.Interior.Color = 13998939

$rng.FormatConditions.Item(3).Delete()
$rng.FormatConditions.Item(2).Delete()

Write-Output "edit complete"
